# Add ALU code fix for gate operations
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# and (row 12): func field (5:0) was 11000, should be 100100
$ws.Cells.Item(12, 7).Value = 100100

# or (row 14): func field (5:0) was 11001, should be 100101
$ws.Cells.Item(14, 7).Value = 100101

# ori (row 15): op code (31:26) was 1110, should be 1101
$ws.Cells.Item(15, 2).Value = 1101

# xor (row 16): func field (5:0) was 11010, should be 100110
$ws.Cells.Item(16, 7).Value = 100110

# xori (row 17): op code (31:26) was 111110, should be 1110
$ws.Cells.Item(17, 2).Value = 1110

# nor (row 18): func field (5:0) was 11011, should be 100111
$ws.Cells.Item(18, 7).Value = 100111

# row 30 func/value fix
$ws.Cells.Item(30, 7).Value = 1001

# Row 14 picked up a manual row-height tweak in the edit session
$ws.Rows.Item(14).RowHeight = 17

# View state: zoom + new selection on the gate-operations area
$excel.ActiveWindow.Zoom = 130
$ws.Range("H5").Select()
